# Rename the worksheets to match the updated naming convention used by the
# refreshed "basic model data" extraction logic.
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("ColumnNames").Name = "column names"
$wb.Worksheets.Item("DerivedData").Name = "derived data"
